$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: duplicate the existing row 2 (mariadb-2-vm) data down into the
#     new row 3, since the sheet is growing from 1 data row to 2 data rows.
$ws.Range("A3").Value = "mariadb-2-vm"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = "Google Compute Engine"
$ws.Range("D3").Value = "mariadb-2-vm"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 21.49573552293964
$ws.Range("K3").Value = 10.28462705365895
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 0.6997745921401075
$ws.Range("N3").Value = 0.6694750599138971
$ws.Range("P3").Value = "Production"
$ws.Range("Q3").Value = 40

# --- Step 2: overwrite row 2 with the newly discovered instance, clearing
#     out the cpu/ram utilization figures that are not yet known for it and
#     recording its disk size/iops related fields.
$ws.Range("A2").Value = "instance-20240725-115904"
$ws.Range("D2").Value = "instance-20240725-115904"
$ws.Range("E2").Value = "X86_64"
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").Value = 1
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = 40
